# MercuryTourFlights.xlsx edit:
#   - Insert a new header row at the top of Sheet1 (from / fromM / fromD /
#     return / returnM / returnD), shifting the existing flight rows down.
#   - Append a new data row (row 4) mirroring the first flight
#     (New York/July -> London/July) but with the day values stored as the
#     text strings "2" and "7" instead of numbers.
#   - Update the active selection on Sheet1 to C18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row 1 for headers; existing rows 1-3 shift down to 2-4.
$ws.Rows.Item(1).Insert()

$ws.Cells.Item(1, 1).Value = "from"
$ws.Cells.Item(1, 2).Value = "fromM"
$ws.Cells.Item(1, 3).Value = "fromD"
$ws.Cells.Item(1, 4).Value = "return"
$ws.Cells.Item(1, 5).Value = "returnM"
$ws.Cells.Item(1, 6).Value = "returnD"

# New row 4: same trip as row 2 but with text day values "2" / "7".
$ws.Cells.Item(4, 1).Value = "New York"
$ws.Cells.Item(4, 2).Value = "July"
$ws.Cells.Item(4, 3).Value = "2"
$ws.Cells.Item(4, 4).Value = "London"
$ws.Cells.Item(4, 5).Value = "July"
$ws.Cells.Item(4, 6).Value = "7"

# Update saved selection to match the authored workbook.
[void]$ws.Range("C18").Select()

Write-Output "edit applied"
